# New crime data collected
# Updates the weekly CompStat report: bumps the report volume/week dates
# in the title block, and refreshes the crime-count / percent-change grid
# (rows 14-30) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Title block: "Volume 30   Number  5" -> "...6" and the reporting
#    week dates "1/30/2023 .. 2/5/2023" -> "2/6/2023 .. 2/12/2023".
#    These live as runs inside a shared string alongside other text, so
#    we patch just the relevant characters in place rather than
#    overwriting the whole cell (which would also nuke the rest of the
#    sentence).
# ---------------------------------------------------------------------
$a8 = $ws.Range("A8")
$volText = $a8.Text
$volPos = $volText.LastIndexOf("5") + 1
$a8.Characters($volPos, 1).Text = "6"

$c9 = $ws.Range("C9")
$weekText = $c9.Text
$endPos = $weekText.IndexOf("2/5/2023") + 1
$c9.Characters($endPos, 8).Text = "2/12/2023"
$weekText2 = $c9.Text
$startPos = $weekText2.IndexOf("1/30/2023") + 1
$c9.Characters($startPos, 9).Text = "2/6/2023"

# ---------------------------------------------------------------------
# 2) Crime grid (rows 14-30). Most cells are plain same-format value
#    swaps; a handful toggle between the "no data" placeholder text
#    ("0" / "***.*", general-formatted) and a real formatted number,
#    which needs the number format flipped before the write so the new
#    value keeps the right stored type.
# ---------------------------------------------------------------------

$numericUpdates = @{
    "G15" = 4;   "H15" = -75;  "I15" = 3;   "J15" = 6;   "K15" = -50;
    "L15" = 200; "M15" = 50;   "N15" = 0;
    "F16" = 9;   "G16" = 7;    "H16" = 28.571428571428;
    "I16" = 13;  "J16" = 13;
    "L16" = 44.444444444444;  "M16" = -40.909090909090; "N16" = -88.181818181818;
    "C17" = 5;   "D17" = 3;    "E17" = 66.666666666666;
    "F17" = 19;  "G17" = 10;   "H17" = 90;
    "I17" = 24;  "J17" = 13;   "K17" = 84.615384615384;
    "L17" = 9.090909090909;   "M17" = 200; "N17" = 0;
    "C18" = 8;   "D18" = 5;    "E18" = 60;
    "F18" = 26;  "G18" = 17;   "H18" = 52.941176470588;
    "I18" = 30;  "J18" = 20;   "K18" = 50;
    "L18" = 150; "M18" = -38.775510204081; "N18" = -83.240223463687;
    "C19" = 14;  "D19" = 13;   "E19" = 7.692307692307;
    "F19" = 44;  "G19" = 54;   "H19" = -18.518518518518;
    "I19" = 65;  "J19" = 100;  "K19" = -35;
    "L19" = 32.653061224489; "M19" = 12.068965517241; "N19" = -23.529411764705;
    "F20" = 9;   "G20" = 5;    "H20" = 80;
    "M20" = -42.105263157894; "N20" = -94.685990338164;
    "C21" = 29;  "E21" = 16;
    "F21" = 108; "G21" = 97;   "H21" = 11.340206185567;
    "I21" = 146; "J21" = 162;  "K21" = -9.876543209876;
    "L21" = 52.083333333333; "M21" = -10.429447852760; "N21" = -76.026272577996;
    "J22" = 4;   "K22" = -75;
    "C24" = 28;  "D24" = 41;   "E24" = -31.707317073170;
    "F24" = 127; "G24" = 140;  "H24" = -9.285714285714;
    "I24" = 179; "J24" = 202;  "K24" = -11.386138613861;
    "L24" = 62.727272727272; "M24" = 72.115384615384;
    "C25" = 10;  "D25" = 9;    "E25" = 11.111111111111;
    "F25" = 34;  "G25" = 39;   "H25" = -12.820512820512;
    "I25" = 52;  "J25" = 54;   "K25" = -3.703703703703;
    "L25" = 18.181818181818; "M25" = 8.333333333333;
    "G26" = 4;   "H26" = -75;  "I26" = 3;   "J26" = 6;   "K26" = -50;
    "L26" = 200;
    "F27" = 1;   "G27" = 1;    "H27" = 0;
    "G30" = 2;   "J30" = 3;
}

foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}

# Cells that were the "N/A" placeholder (text) and become real counts /
# percent-changes. Give them the right number format first so the write
# is stored as a number, not re-parsed back into placeholder text.
$countFmt = "#,##0"
$pctFmt = "#,##0.0;""-""#,##0.0"

$toNumberCount = @{
    "C15" = 1; "D15" = 3;
    "D22" = 1;
    "C26" = 1; "D26" = 3;
    "D30" = 2;
}
foreach ($ref in $toNumberCount.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = $countFmt
    $c.Value = $toNumberCount[$ref]
}

$toNumberPct = @{
    "N14" = -100;
    "E15" = -66.666666666666;
    "E22" = -100;
    "E26" = -66.666666666666;
    "E30" = -100;
}
foreach ($ref in $toNumberPct.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = $pctFmt
    $c.Value = $toNumberPct[$ref]
}

# Cells that go the other way: a real number reverts to the "no data"
# placeholder text. Force text storage via NumberFormat "@" and then
# restore the original (general) look by copying the format from a
# known-good placeholder cell (A14) so the style index matches the
# other untouched placeholder cells.
$toPlaceholder = @{
    "C20" = "0";
    "D20" = "0";
    "E20" = "***.*";
    "C28" = "0";
    "C29" = "0";
}
$styleDonor = $ws.Range("A14")
foreach ($ref in $toPlaceholder.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $toPlaceholder[$ref]
    $styleDonor.Copy() | Out-Null
    $c.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false
